$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 6071.7144
$ws.Cells.Item(70, 9).Value = 4002
$ws.Cells.Item(70, 10).Value = 6416.6665
$ws.Cells.Item(70, 11).Value = 12006
$ws.Cells.Item(70, 12).Value = 19249.9995
$ws.Cells.Item(70, 13).Value = -11736
$ws.Cells.Item(70, 14).Value = -19789.9995
$ws.Cells.Item(73, 8).Value = 6071.7144
$ws.Cells.Item(73, 9).Value = 4002
$ws.Cells.Item(73, 10).Value = 6416.6665
$ws.Cells.Item(73, 11).Value = 12006
$ws.Cells.Item(73, 12).Value = 19249.9995
$ws.Cells.Item(73, 13).Value = -11070
$ws.Cells.Item(73, 14).Value = -21121.9995
$ws.Cells.Item(127, 8).Value = 2503.8667
$ws.Cells.Item(127, 9).Value = 2327.0908
$ws.Cells.Item(127, 11).Value = 6981.2724
$ws.Cells.Item(127, 13).Value = -2021.2724
$ws.Cells.Item(137, 8).Value = 11499708
$ws.Cells.Item(137, 9).Value = 2232.8
$ws.Cells.Item(137, 10).Value = 13895015
$ws.Cells.Item(137, 11).Value = 6698.400000000001
$ws.Cells.Item(137, 12).Value = 41685045
$ws.Cells.Item(137, 13).Value = -4148.400000000001
$ws.Cells.Item(137, 14).Value = -41690145

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5921.851
$ws.Cells.Item(32, 9).Value = 3884.606
$ws.Cells.Item(32, 11).Value = 3884.606
$ws.Cells.Item(32, 13).Value = -3597.606
$ws.Cells.Item(61, 8).Value = 39352.75
$ws.Cells.Item(61, 9).Value = 39352.75
$ws.Cells.Item(61, 11).Value = 39352.75
$ws.Cells.Item(61, 13).Value = -39140.75
$ws.Cells.Item(75, 8).Value = 20000
$ws.Cells.Item(75, 10).Value = 20000
$ws.Cells.Item(75, 12).Value = 20000
$ws.Cells.Item(75, 14).Value = -21748
$ws.Cells.Item(78, 8).Value = 20000
$ws.Cells.Item(78, 10).Value = 20000
$ws.Cells.Item(78, 12).Value = 60000
$ws.Cells.Item(78, 14).Value = -68736
$ws.Cells.Item(122, 8).Value = 3593.875
$ws.Cells.Item(122, 9).Value = 2437.8
$ws.Cells.Item(122, 11).Value = 7313.400000000001
$ws.Cells.Item(122, 13).Value = -4863.400000000001
$ws.Cells.Item(132, 8).Value = 3174.95
$ws.Cells.Item(132, 9).Value = 1785.3
$ws.Cells.Item(132, 10).Value = 4564.6
$ws.Cells.Item(132, 11).Value = 5355.9
$ws.Cells.Item(132, 12).Value = 13693.8
$ws.Cells.Item(132, 13).Value = -2825.9
$ws.Cells.Item(132, 14).Value = -18753.8
$ws.Cells.Item(135, 8).Value = 86919.336
$ws.Cells.Item(135, 10).Value = 86919.336
$ws.Cells.Item(135, 12).Value = 86919.336
$ws.Cells.Item(135, 14).Value = -97059.336
$ws.Cells.Item(136, 8).Value = 39352.75
$ws.Cells.Item(136, 9).Value = 39352.75
$ws.Cells.Item(136, 11).Value = 118058.25
$ws.Cells.Item(136, 13).Value = -115508.25
$ws.Cells.Item(139, 8).Value = 85683.25
$ws.Cells.Item(139, 10).Value = 72928.28999999999
$ws.Cells.Item(139, 12).Value = 72928.28999999999
$ws.Cells.Item(139, 14).Value = -83208.28999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(81, 8).Value = 44061.445
$ws.Cells.Item(81, 10).Value = 44061.445
$ws.Cells.Item(81, 12).Value = 44061.445
$ws.Cells.Item(81, 14).Value = -46183.445
$ws.Cells.Item(84, 8).Value = 44061.445
$ws.Cells.Item(84, 10).Value = 44061.445
$ws.Cells.Item(84, 12).Value = 132184.335
$ws.Cells.Item(84, 14).Value = -142792.335
$ws.Cells.Item(105, 8).Value = 2755.889
$ws.Cells.Item(105, 9).Value = 2694.125
$ws.Cells.Item(105, 11).Value = 2694.125
$ws.Cells.Item(105, 13).Value = -947.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2233.7778
$ws.Cells.Item(16, 9).Value = 2096.375
$ws.Cells.Item(16, 11).Value = 2096.375
$ws.Cells.Item(16, 13).Value = -1809.375
$ws.Cells.Item(31, 8).Value = 28578514
$ws.Cells.Item(31, 9).Value = 142858240
$ws.Cells.Item(31, 10).Value = 8582.429
$ws.Cells.Item(31, 11).Value = 142858240
$ws.Cells.Item(31, 12).Value = 8582.429
$ws.Cells.Item(31, 13).Value = -142857945
$ws.Cells.Item(31, 14).Value = -9172.429
$ws.Cells.Item(34, 8).Value = 28578514
$ws.Cells.Item(34, 9).Value = 142858240
$ws.Cells.Item(34, 10).Value = 8582.429
$ws.Cells.Item(34, 11).Value = 142858240
$ws.Cells.Item(34, 12).Value = 8582.429
$ws.Cells.Item(34, 13).Value = -142858038
$ws.Cells.Item(34, 14).Value = -8986.429
$ws.Cells.Item(62, 8).Value = 23804.133
$ws.Cells.Item(62, 10).Value = 31108
$ws.Cells.Item(62, 12).Value = 31108
$ws.Cells.Item(62, 14).Value = -32356
$ws.Cells.Item(65, 8).Value = 23804.133
$ws.Cells.Item(65, 10).Value = 31108
$ws.Cells.Item(65, 12).Value = 155540
$ws.Cells.Item(65, 14).Value = -161780
$ws.Cells.Item(99, 8).Value = 8777.777
$ws.Cells.Item(99, 9).Value = 4333.3335
$ws.Cells.Item(99, 11).Value = 4333.3335
$ws.Cells.Item(99, 13).Value = -2835.3335
$ws.Cells.Item(113, 8).Value = 2233.7778
$ws.Cells.Item(113, 9).Value = 2096.375
$ws.Cells.Item(113, 11).Value = 2096.375
$ws.Cells.Item(113, 13).Value = 73.625
$ws.Cells.Item(122, 8).Value = 3237.375
$ws.Cells.Item(122, 9).Value = 2032.1428
$ws.Cells.Item(122, 11).Value = 6096.428400000001
$ws.Cells.Item(122, 13).Value = -3646.428400000001
$ws.Cells.Item(126, 8).Value = 8777.777
$ws.Cells.Item(126, 9).Value = 4333.3335
$ws.Cells.Item(126, 11).Value = 13000.0005
$ws.Cells.Item(126, 13).Value = -10530.0005
$ws.Cells.Item(134, 8).Value = 4178.8887
$ws.Cells.Item(134, 9).Value = 4363.75
$ws.Cells.Item(134, 11).Value = 13091.25
$ws.Cells.Item(134, 13).Value = -10556.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(25, 8).Value = 2469.4
$ws.Cells.Item(25, 9).Value = 150
$ws.Cells.Item(25, 10).Value = 3049.25
$ws.Cells.Item(25, 11).Value = 450
$ws.Cells.Item(25, 12).Value = 9147.75
$ws.Cells.Item(25, 13).Value = -281
$ws.Cells.Item(25, 14).Value = -9485.75
$ws.Cells.Item(30, 8).Value = 2469.4
$ws.Cells.Item(30, 9).Value = 150
$ws.Cells.Item(30, 10).Value = 3049.25
$ws.Cells.Item(30, 11).Value = 450
$ws.Cells.Item(30, 12).Value = 9147.75
$ws.Cells.Item(30, 13).Value = -348
$ws.Cells.Item(30, 14).Value = -9351.75
$ws.Cells.Item(80, 8).Value = 5811.75
$ws.Cells.Item(80, 9).Value = 5499
$ws.Cells.Item(80, 10).Value = 5874.3
$ws.Cells.Item(80, 11).Value = 16497
$ws.Cells.Item(80, 12).Value = 17622.9
$ws.Cells.Item(80, 13).Value = -15561
$ws.Cells.Item(80, 14).Value = -19494.9
$ws.Cells.Item(83, 8).Value = 5811.75
$ws.Cells.Item(83, 9).Value = 5499
$ws.Cells.Item(83, 10).Value = 5874.3
$ws.Cells.Item(83, 11).Value = 49491
$ws.Cells.Item(83, 12).Value = 52868.7
$ws.Cells.Item(83, 13).Value = -44811
$ws.Cells.Item(83, 14).Value = -62228.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 13).Value = ""
$ws.Cells.Item(70, 8).Value = 1836610.2
$ws.Cells.Item(70, 9).Value = 2805992.8
$ws.Cells.Item(70, 10).Value = 5554.1113
$ws.Cells.Item(70, 11).Value = 2805992.8
$ws.Cells.Item(70, 12).Value = 5554.1113
$ws.Cells.Item(70, 13).Value = -2805722.8
$ws.Cells.Item(70, 14).Value = -6094.1113
$ws.Cells.Item(73, 8).Value = 1836610.2
$ws.Cells.Item(73, 9).Value = 2805992.8
$ws.Cells.Item(73, 10).Value = 5554.1113
$ws.Cells.Item(73, 11).Value = 2805992.8
$ws.Cells.Item(73, 12).Value = 5554.1113
$ws.Cells.Item(73, 13).Value = -2805056.8
$ws.Cells.Item(73, 14).Value = -7426.1113
$ws.Cells.Item(102, 8).Value = 20841096
$ws.Cells.Item(102, 9).Value = 31258506
$ws.Cells.Item(102, 11).Value = 31258506
$ws.Cells.Item(102, 13).Value = -31256884
$ws.Cells.Item(122, 8).Value = 410417.28
$ws.Cells.Item(122, 9).Value = 526465.4399999999
$ws.Cells.Item(122, 10).Value = 4248.8335
$ws.Cells.Item(122, 11).Value = 1579396.32
$ws.Cells.Item(122, 12).Value = 12746.5005
$ws.Cells.Item(122, 13).Value = -1576946.32
$ws.Cells.Item(122, 14).Value = -17646.5005
$ws.Cells.Item(126, 8).Value = 4929.1904
$ws.Cells.Item(126, 9).Value = 2803.6428
$ws.Cells.Item(126, 10).Value = 9180.286
$ws.Cells.Item(126, 11).Value = 8410.928400000001
$ws.Cells.Item(126, 12).Value = 27540.858
$ws.Cells.Item(126, 13).Value = -5940.928400000001
$ws.Cells.Item(126, 14).Value = -32480.858
$ws.Cells.Item(132, 8).Value = 3298.3235
$ws.Cells.Item(132, 9).Value = 2676.25
$ws.Cells.Item(132, 11).Value = 8028.75
$ws.Cells.Item(132, 13).Value = -5498.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5321.5
$ws.Cells.Item(7, 9).Value = 1999.5
$ws.Cells.Item(7, 11).Value = 1999.5
$ws.Cells.Item(7, 13).Value = -1887.5
$ws.Cells.Item(126, 8).Value = 5321.5
$ws.Cells.Item(126, 9).Value = 1999.5
$ws.Cells.Item(126, 11).Value = 5998.5
$ws.Cells.Item(126, 13).Value = -3528.5
$ws.Cells.Item(136, 8).Value = 7590.727
$ws.Cells.Item(136, 9).Value = 2500
$ws.Cells.Item(136, 11).Value = 7500
$ws.Cells.Item(136, 13).Value = -4950

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 6833.5
$ws.Cells.Item(2, 9).Value = 7000.1816
$ws.Cells.Item(2, 11).Value = 7000.1816
$ws.Cells.Item(2, 13).Value = -6888.1816
$ws.Cells.Item(81, 8).Value = 5211883
$ws.Cells.Item(81, 9).Value = 4167546.5
$ws.Cells.Item(81, 10).Value = 6952443.5
$ws.Cells.Item(81, 11).Value = 8335093
$ws.Cells.Item(81, 12).Value = 13904887
$ws.Cells.Item(81, 13).Value = -8334032
$ws.Cells.Item(81, 14).Value = -13907009
$ws.Cells.Item(84, 8).Value = 5211883
$ws.Cells.Item(84, 9).Value = 4167546.5
$ws.Cells.Item(84, 10).Value = 6952443.5
$ws.Cells.Item(84, 11).Value = 41675465
$ws.Cells.Item(84, 12).Value = 69524435
$ws.Cells.Item(84, 13).Value = -41670161
$ws.Cells.Item(84, 14).Value = -69535043
